$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2309.375
$ws.Range("I40").Value = 1906.25
$ws.Range("K40").Value = 1906.25
$ws.Range("M40").Value = -1731.25
$ws.Range("H70").Value = 185993
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("M70").Value = -8730
$ws.Range("H73").Value = 185993
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("M73").Value = -8064
$ws.Range("H80").Value = 1094.3914
$ws.Range("I80").Value = 1139.5294
$ws.Range("K80").Value = 3418.5882
$ws.Range("M80").Value = -2420.5882
$ws.Range("H83").Value = 1094.3914
$ws.Range("I83").Value = 1139.5294
$ws.Range("K83").Value = 10255.7646
$ws.Range("M83").Value = -5263.764599999999
$ws.Range("H86").Value = 2846.8333
$ws.Range("I86").Value = 2216.2
$ws.Range("K86").Value = 2216.2
$ws.Range("M86").Value = -1093.2
$ws.Range("H89").Value = 2846.8333
$ws.Range("I89").Value = 2216.2
$ws.Range("K89").Value = 11081
$ws.Range("M89").Value = -5465
$ws.Range("H103").Value = 1168
$ws.Range("I103").Value = 1168
$ws.Range("K103").Value = 3504
$ws.Range("M103").Value = -2918
$ws.Range("H137").Value = 1563.2
$ws.Range("I137").Value = 1446.6
$ws.Range("J137").Value = 1679.8
$ws.Range("K137").Value = 4339.799999999999
$ws.Range("L137").Value = 5039.4
$ws.Range("M137").Value = -1789.799999999999
$ws.Range("N137").Value = -10139.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1027.5
$ws.Range("I2").Value = 1027.5
$ws.Range("K2").Value = 1027.5
$ws.Range("M2").Value = -914.5
$ws.Range("H45").Value = 4437.375
$ws.Range("I45").Value = 4437.375
$ws.Range("K45").Value = 4437.375
$ws.Range("M45").Value = -4060.375
$ws.Range("H97").Value = 1722.25
$ws.Range("I97").Value = 2136.6667
$ws.Range("K97").Value = 2136.6667
$ws.Range("M97").Value = -1640.6667
$ws.Range("H110").Value = 12222.25
$ws.Range("I110").Value = 16449.5
$ws.Range("K110").Value = 16449.5
$ws.Range("M110").Value = -14404.5
$ws.Range("H116").Value = 1027.5
$ws.Range("I116").Value = 1027.5
$ws.Range("K116").Value = 1027.5
$ws.Range("M116").Value = 1266.5
$ws.Range("H122").Value = 502434.7
$ws.Range("I122").Value = 590570.25
$ws.Range("K122").Value = 1771710.75
$ws.Range("M122").Value = -1769260.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1027.5
$ws.Range("I3").Value = 1027.5
$ws.Range("K3").Value = 1027.5
$ws.Range("M3").Value = -913.5
$ws.Range("H86").Value = 2285.4285
$ws.Range("I86").Value = 2249.25
$ws.Range("J86").Value = 2333.6667
$ws.Range("K86").Value = 2249.25
$ws.Range("L86").Value = 2333.6667
$ws.Range("M86").Value = -1126.25
$ws.Range("N86").Value = -4579.6667
$ws.Range("H89").Value = 2285.4285
$ws.Range("I89").Value = 2249.25
$ws.Range("J89").Value = 2333.6667
$ws.Range("K89").Value = 11246.25
$ws.Range("L89").Value = 11668.3335
$ws.Range("M89").Value = -5630.25
$ws.Range("N89").Value = -22900.3335

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4426.7896
$ws.Range("I31").Value = 3964.2144
$ws.Range("K31").Value = 3964.2144
$ws.Range("M31").Value = -3669.2144
$ws.Range("H34").Value = 4426.7896
$ws.Range("I34").Value = 3964.2144
$ws.Range("K34").Value = 3964.2144
$ws.Range("M34").Value = -3762.2144
$ws.Range("H62").Value = 52901.125
$ws.Range("I62").Value = 3127.5
$ws.Range("J62").Value = 102674.75
$ws.Range("K62").Value = 3127.5
$ws.Range("L62").Value = 102674.75
$ws.Range("M62").Value = -2503.5
$ws.Range("N62").Value = -103922.75
$ws.Range("H65").Value = 52901.125
$ws.Range("I65").Value = 3127.5
$ws.Range("J65").Value = 102674.75
$ws.Range("K65").Value = 15637.5
$ws.Range("L65").Value = 513373.75
$ws.Range("M65").Value = -12517.5
$ws.Range("N65").Value = -519613.75
$ws.Range("H68").Value = 46000
$ws.Range("J68").Value = 46000
$ws.Range("L68").Value = 46000
$ws.Range("N68").Value = -47498
$ws.Range("H71").Value = 46000
$ws.Range("J71").Value = 46000
$ws.Range("L71").Value = 138000
$ws.Range("N71").Value = -145488
$ws.Range("H99").Value = 13555.044
$ws.Range("I99").Value = 9978.299999999999
$ws.Range("K99").Value = 9978.299999999999
$ws.Range("M99").Value = -8480.299999999999
$ws.Range("H107").Value = 1181
$ws.Range("I107").Value = 738.75
$ws.Range("K107").Value = 738.75
$ws.Range("M107").Value = 1181.25
$ws.Range("H122").Value = 3652.6365
$ws.Range("I122").Value = 3767.9
$ws.Range("K122").Value = 11303.7
$ws.Range("M122").Value = -8853.700000000001
$ws.Range("H126").Value = 13555.044
$ws.Range("I126").Value = 9978.299999999999
$ws.Range("K126").Value = 29934.9
$ws.Range("M126").Value = -27464.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47810.715
$ws.Range("I2").Value = 71466.42999999999
$ws.Range("J2").Value = 499.2857
$ws.Range("K2").Value = 428798.58
$ws.Range("L2").Value = 2995.7142
$ws.Range("M2").Value = -428685.58
$ws.Range("N2").Value = -3221.7142
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H4").Value = 636796.4399999999
$ws.Range("I4").Value = 954790.9
$ws.Range("J4").Value = 807.63635
$ws.Range("K4").Value = 2864372.7
$ws.Range("L4").Value = 2422.90905
$ws.Range("M4").Value = -2864260.7
$ws.Range("N4").Value = -2646.90905
$ws.Range("H12").Value = 257.60605
$ws.Range("I12").Value = 217.91667
$ws.Range("J12").Value = 280.2857
$ws.Range("K12").Value = 653.75001
$ws.Range("L12").Value = 840.8571000000001
$ws.Range("M12").Value = -480.75001
$ws.Range("N12").Value = -1186.8571
$ws.Range("H14").Value = 911.3
$ws.Range("I14").Value = 911.3
$ws.Range("K14").Value = 2733.9
$ws.Range("M14").Value = -2560.9
$ws.Range("H133").Value = 1200
$ws.Range("I133").Value = 1200
$ws.Range("K133").Value = 3600
$ws.Range("M133").Value = 1460

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7214.2856
$ws.Range("I70").Value = 5250
$ws.Range("K70").Value = 5250
$ws.Range("M70").Value = -4980
$ws.Range("H73").Value = 7214.2856
$ws.Range("I73").Value = 5250
$ws.Range("K73").Value = 5250
$ws.Range("M73").Value = -4314
$ws.Range("H97").Value = 402.66666
$ws.Range("J97").Value = 402.66666
$ws.Range("L97").Value = 402.66666
$ws.Range("N97").Value = -1394.66666
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
$ws.Range("H122").Value = 35088.645
$ws.Range("I122").Value = 2574.3333
$ws.Range("J122").Value = 146566.28
$ws.Range("K122").Value = 7722.999899999999
$ws.Range("L122").Value = 439698.84
$ws.Range("M122").Value = -5272.999899999999
$ws.Range("N122").Value = -444598.84
$ws.Range("H123").Value = 28885.555
$ws.Range("J123").Value = 28885.555
$ws.Range("L123").Value = 28885.555
$ws.Range("N123").Value = -33785.555

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5079.8
$ws.Range("I7").Value = 5079.8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5079.8
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4967.8
$ws.Range("H46").Value = 4379.4443
$ws.Range("I46").Value = 1998.75
$ws.Range("J46").Value = 6284
$ws.Range("K46").Value = 1998.75
$ws.Range("L46").Value = 6284
$ws.Range("M46").Value = -1810.75
$ws.Range("N46").Value = -6660
$ws.Range("H61").Value = 3796.375
$ws.Range("I61").Value = 3794.5
$ws.Range("J61").Value = 3802
$ws.Range("K61").Value = 3794.5
$ws.Range("L61").Value = 3802
$ws.Range("M61").Value = -3592.5
$ws.Range("N61").Value = -4206
$ws.Range("H93").Value = 1346.4615
$ws.Range("I93").Value = 1230.4
$ws.Range("K93").Value = 1230.4
$ws.Range("M93").Value = 17.59999999999991
$ws.Range("H113").Value = 3796.375
$ws.Range("I113").Value = 3794.5
$ws.Range("J113").Value = 3802
$ws.Range("K113").Value = 3794.5
$ws.Range("L113").Value = 3802
$ws.Range("M113").Value = -1624.5
$ws.Range("N113").Value = -8142
$ws.Range("H126").Value = 5079.8
$ws.Range("I126").Value = 5079.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15239.4
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -12769.4
$ws.Range("H132").Value = 3966.6667
$ws.Range("I132").Value = 3560
$ws.Range("K132").Value = 10680
$ws.Range("M132").Value = -8150

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1329.875
$ws.Range("I107").Value = 1278.6
$ws.Range("J107").Value = 1415.3334
$ws.Range("K107").Value = 3835.8
$ws.Range("L107").Value = 4246.0002
$ws.Range("M107").Value = -1915.8
$ws.Range("N107").Value = -8086.0002
$ws.Range("H113").Value = 561.05
$ws.Range("I113").Value = 613.1111
$ws.Range("K113").Value = 1839.3333
$ws.Range("M113").Value = 330.6667000000002
